$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New gradebook columns: "Section 7" (M) and "MidTerm2" (N) ---
$ws.Cells.Item(1, 13).Value = "Section 7"
$ws.Cells.Item(1, 14).Value = "MidTerm2"

$newScores = @{
    2  = @(2, 53.5)
    3  = @(2, 65.5)
    4  = @(2, 70.5)
    5  = @(2, 62.5)
    6  = @(2, 68)
    7  = @(2, 73)
    8  = @(2, 61.5)
    9  = @(2, 71.5)
    10 = @(0, 64)
    11 = @(2, 60.5)
    12 = @(2, 65.5)
    13 = @(2, 62)
    14 = @(2, 69)
    15 = @(2, 64)
    16 = @(2, 65)
    17 = @(2, 68)
}

foreach ($row in $newScores.Keys) {
    $vals = $newScores[$row]
    $ws.Cells.Item($row, 13).Value = $vals[0]
    $ws.Cells.Item($row, 14).Value = $vals[1]
}

# --- Un-hide the previously-hidden section columns (C:K) and give them a
#     normal, visible width now that the gradebook is being reviewed ---
for ($c = 3; $c -le 11; $c++) {
    $col = $ws.Columns.Item($c)
    $col.Hidden = $false
    $col.ColumnWidth = 10
}

# --- Update the saved selection/view state ---
[void]$ws.Range("N18").Select()
